$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.158716678619385
$ws.Range("B1").Value = 2.369450569152832
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.39897894859314
$ws.Range("E1").Value = 1.221897125244141
